$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: text code changes from "004" to "001".
# Temporarily force a text number-format so the numeric-looking string is
# kept as text (matches the source data, which stores it as a text value),
# then clear the formatting again so the cell is left exactly as before.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# M2 / N2: date-like text fields (stored as plain text, not real dates)
$ws.Range("M2").Value = "2020-12-15 00:00:00"
$ws.Range("N2").Value = "2018-12-31 00:00:00"

# O2:AD2 numeric figures updated to the new reported values
$ws.Range("O2").Value = 265948324.95
$ws.Range("P2").Value = 299.4623182198
$ws.Range("Q2").Value = 1979181334.04
$ws.Range("R2").Value = 2228.5917032206
$ws.Range("S2").Value = 292848187.77
$ws.Range("T2").Value = 329.752019354
$ws.Range("U2").Value = -138977008.98
$ws.Range("V2").Value = -156.4904659438
$ws.Range("W2").Value = 1917.81
$ws.Range("X2").Value = 0.0021594865
$ws.Range("Y2").Value = 148433882.62
$ws.Range("Z2").Value = 167.1390658321
$ws.Range("AA2").Value = -38405587.74
$ws.Range("AB2").Value = -43.2453422648
$ws.Range("AC2").Value = 88808610.89
$ws.Range("AD2").Value = 57.2754815426
